$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Tabelle1"

# Row 1 - AI-guided web crawler thesis
$ws.Range("A1").Value = "AI-GUIDED WEB CRAWLER FOR AUTOMATIC DETECTION OF MALICIOUS SITES"
$ws.Range("B1").Value = "d279620"
$ws.Range("C1").Value = "research project"
$ws.Range("D1").Value = "This thesis focuses on developing an AI-guided web crawler for the automatic detection of malicious sites. The research aims to leverage artificial intelligence to enhance the efficiency and accuracy of web crawling in identifying and cataloging potentially harmful websites."
$ws.Range("E1").Value = "web development, cybersecurity, and machine learning"
$ws.Range("F1").Value = "The project involves implementing machine learning algorithms for pattern recognition, collaborating with cybersecurity experts, and optimizing web crawling algorithms for real-time detection"
$ws.Range("G1").Value = (Get-Date -Year 2024 -Month 10 -Day 11 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("H1").Value = "LM"
$ws.Range("I1").Value = "LM-32"

# Row 2 - Kafka clients reactive API
$ws.Range("A2").Value = "PERFORMANCE EVALUATION OF KAFKA CLIENTS USING A REACTIVE API"
$ws.Range("B2").Value = "d370335"
$ws.Range("C2").Value = "research project"
$ws.Range("D2").Value = "This thesis focuses on the performance evaluation of Kafka clients using a reactive API. The research aims to assess and enhance the efficiency of Kafka clients by implementing a reactive programming approach. The study explores how a reactive API can improve responsiveness and scalability in real-time data streaming applications."
$ws.Range("E2").Value = "networking protocols, congestion control algorithms, and familiarity with QUIC"
$ws.Range("F2").Value = "The study involves simulations, performance evaluations, and an in-depth analysis of the effectiveness of different congestion control schemes in QUIC"
$ws.Range("G1").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = (Get-Date -Year 2025 -Month 12 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("H2").Value = "LM"
$ws.Range("I2").Value = "LM-32"

# Row 3 - Congestion control schemes on QUIC
$ws.Range("A3").Value = "A STUDY OF CONGESTION CONTROL SCHEMES ON QUIC"
$ws.Range("B3").Value = "d350985"
$ws.Range("C3").Value = "research project"
$ws.Range("D3").Value = "This research paper delves into a comprehensive study of congestion control schemes on QUIC (Quick UDP Internet Connections). The investigation aims to analyze and compare various congestion control strategies within the QUIC protocol, shedding light on their impact on network performance and efficiency."
$ws.Range("G1").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = (Get-Date -Year 2026 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("H3").Value = "LM"
$ws.Range("I3").Value = "LM-32"

# Row 4 - Ultra-low-power acoustic sensor frontend
$ws.Range("A4").Value = "ULTRA-LOW-POWER ACOUSTIC SENSOR FRONTEND A DIGITAL TRANSCONDUCTANCE AMPLIFIER APPROACH"
$ws.Range("B4").Value = "d255269"
$ws.Range("C4").Value = "dissertation"
$ws.Range("D4").Value = "This dissertation focuses on the design and implementation of an ultra-low-power acoustic sensor frontend using a digital transconductance amplifier approach. The research aims to explore novel methods in signal processing and circuit design to achieve highly efficient and low-power acoustic sensing capabilities."
$ws.Range("E4").Value = "analog and digital circuit design, signal processing, and low-power electronics"
$ws.Range("F4").Value = "The project involves the development of a digital transconductance amplifier, integration with an acoustic sensor frontend, and extensive testing for ultra-low-power performance"
$ws.Range("G1").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G4").Value = (Get-Date -Year 2027 -Month 9 -Day 29 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("H4").Value = "LM"
$ws.Range("I4").Value = "LM-31"

# Row 5 - Preliminary design of an arcjet
$ws.Range("A5").Value = "PRELIMINARY DESIGN OF AN ARCJET IN THE 1KW CLASS FOR SPACE APPLICATION"
$ws.Range("B5").Value = "d357587"
$ws.Range("C5").Value = "engineering project"
$ws.Range("D5").Value = "This engineering project involves the preliminary design of an arcjet in the 1kW class for space applications. The research aims to conceptualize and outline the key parameters, components, and specifications required for the development of a high-powered arcjet propulsion system suitable for space missions."
$ws.Range("E5").Value = "propulsion systems, plasma physics, and aerospace engineering"
$ws.Range("F5").Value = "The project includes conceptual design, performance modeling, and consideration of thermal management aspects for the arcjet propulsion system."
$ws.Range("G1").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Value = (Get-Date -Year 2025 -Month 3 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("H5").Value = "LM"
$ws.Range("I5").Value = "LM-23"

$excel.CutCopyMode = $false
$ws.Range("E14").Select()
